# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a "Label" column (H) that flags each patient row as Control (0) or MDD (1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column, matching the style used by the other headers / row labels
# (bold font, thin border, centered horizontally, top-aligned vertically).
$ws.Range("H1").Value = "Label"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.LineStyle = 1

# Row -> Label value (0 = Control patient, 1 = MDD patient), for both the
# 100-iteration block (rows 2-11) and the 200-iteration block (rows 12-21).
$labels = @{
    2  = 0;  3  = 0;  4  = 0;  5  = 0;  6  = 0
    7  = 1;  8  = 1;  9  = 1;  10 = 1;  11 = 1
    12 = 0;  13 = 0;  14 = 0;  15 = 0;  16 = 0
    17 = 1;  18 = 1;  19 = 1;  20 = 1;  21 = 1
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 8).Value = $labels[$row]
}
